$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bump the "Förändrad" (changed) date in column C for data rows 2-28 by one day
for ($r = 2; $r -le 28; $r++) {
    $ws.Cells.Item($r, 3).Value = 45428
}

# Remove the two trailing rows (30 then 29, so indices stay valid)
$ws.Rows.Item(30).Delete()
$ws.Rows.Item(29).Delete()

# Row 28 is now the last row; drop its explicit custom height so it matches
# the convention used for the final row in the sheet
$ws.Rows.Item(28).AutoFit()
